$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "updated at" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 00:14"

# --- Update country stats (rows correspond to countries in column A) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5600763
$ws.Range("C4").Value = 34709
$ws.Range("D4").Value = 2962903
$ws.Range("E4").Value = 2464256
$ws.Range("G4").Value = 477
$ws.Range("H4").Value = 173604

# Row 5: Brasil
$ws.Range("B5").Value = 3359570
$ws.Range("C5").Value = 19373
$ws.Range("E5").Value = 818578
$ws.Range("G5").Value = 657
$ws.Range("H5").Value = 108536

# Row 22: Arabia Saudita
$ws.Range("B22").Value = 226622
$ws.Range("C22").Value = 1625
$ws.Range("E22").Value = 14426

# Row 23: Argentina
$ws.Range("D23").Value = 84065
$ws.Range("E23").Value = 104535

# Row 53: Republica Dominicana
$ws.Range("B53").Value = 47185
$ws.Range("C53").Value = 350
$ws.Range("D53").Value = 43529
$ws.Range("E53").Value = 3483

# --- Swap "Montserrat" and "Islas Malvinas" rows (213/214) ---
# Before: row213 = Montserrat (D=12,H=1), row214 = Islas Malvinas (D=13,H=0)
# After:  row213 = Islas Malvinas (D=13,H=0), row214 = Montserrat (D=12,H=1)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
